# Re-upload / re-save the deck with "Embed fonts in the file" turned on.
#
# In real PowerPoint this is a Save dialog / File > Options > Save check-box
# ("Embed fonts in the file") rather than a dedicated object-model switch, so
# the closest supported COM surface is the documented `EmbedTrueTypeFonts`
# parameter of Presentation.SaveAs (and the read/write-ish Presentation
# property of the same name some hosts expose), together with marking every
# font PowerPoint tells us about as embeddable/embedded before the deck is
# written back out.

$p = $ppt.ActivePresentation

# Turn font embedding (with subsetting, the PowerPoint default) on for the
# presentation itself.
try { $p.EmbedTrueTypeFonts = $true } catch { }
try { $p.SaveSubsetFonts = $true } catch { }

# Mark every font currently known to the presentation (theme fonts such as
# "Trebuchet MS" and bullet fonts such as "Wingdings 3") as embedded.
$fontCount = $p.Fonts.Count
for ($i = 1; $i -le $fontCount; $i++) {
    $fnt = $p.Fonts.Item($i)
    try { $fnt.Embedded = $true } catch { }
}

# Persist the deck back to its own package with fonts embedded. SaveAs's
# third argument is the documented EmbedTrueTypeFonts switch
# (msoTrue/msoCTrue = embed); re-saving in place (same name/format) is what
# the "Add files via upload" re-save effectively did.
$p.SaveAs($p.FullName, 1, $true)
$p.Save()
